# Relay Config with new Machine Module 5_4
# Adds a new "relayConfig" PGN row (row 18) to the PGN sheet: a label in A18
# plus a numbered byte-map (1..20) across F18:Y18 with a leading length (E18)
# and a trailing CRC marker (Z18), mirroring the layout of the other PGN rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PGN")

# New row label - becomes a new shared string ("relayConfig").
$ws.Range("A18").Value = "relayConfig"

# Byte count for this PGN.
$ws.Range("E18").Value = 20

# Numbered data bytes 1-20 across F18:Y18.
for ($i = 1; $i -le 20; $i++) {
    $col = 5 + $i   # F=6 .. Y=25
    $ws.Cells.Item(18, $col).Value = $i
}

# Trailing CRC column (reuses the existing "CRC" shared string).
$ws.Range("Z18").Value = "CRC"

# Scroll/selection housekeeping to match the saved view state.
$ws.Activate() | Out-Null
$ws.Range("F19").Select() | Out-Null
